# Applies the "bridgeWeight" testing data + b-value Player1-wins data
# to Sheet1, per the commit "Testing for bridge weights and comparing
# MCTS and Minimax".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in missing "Player 1 wins" counts for the existing "b value" table (H column) ---
$ws.Range("H5").Value = 16
$ws.Range("H6").Value = 13
$ws.Range("H7").Value = 14

# --- Column headers (row 4) for the two new bridgeWeight tables ---
$ws.Range("L4").Value = "bridgeWeight"
$ws.Range("M4").Value = "Player 1 wins"
$ws.Range("N4").Value = "Player 2 wins"
$ws.Range("O4").Value = "Player 1 winrate"

$ws.Range("Q4").Value = "bridgeWeight"
$ws.Range("R4").Value = "Player 1 wins"
$ws.Range("S4").Value = "Player 2 wins"
$ws.Range("T4").Value = "Player 1 winrate"

# --- New section headers (row 3): time-control labels above the two new bridgeWeight tables ---
$ws.Range("Q3").Value = "60Sec"
$ws.Range("L3").Value = "10Sec"

# --- bridgeWeight values (L and Q columns), rows 5-13 ---
$bridgeWeights = @(2.1, 2.2, 2.3, 2.4, 2.5, 2.6, 2.7, 2.8, 2.9)
for ($i = 0; $i -lt $bridgeWeights.Length; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 12).Value = $bridgeWeights[$i]   # column L
    $ws.Cells.Item($row, 17).Value = $bridgeWeights[$i]   # column Q
}

# --- "10Sec" table (L:O): Player 1 win counts ---
$p1Wins10Sec = @(16, 16, 20, 15, 16, 20, 16, 14, 10)
for ($i = 0; $i -lt $p1Wins10Sec.Length; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 13).Value = $p1Wins10Sec[$i]     # column M
}
# Row 5 gets its own (non-shared) formula; rows 6-13 fill down together as a shared formula
$ws.Range("N5").Formula = "=30-M5"
$ws.Range("O5").Formula = "=M5/30"
$ws.Range("N6:N13").Formula = "=30-M6"
$ws.Range("O6:O13").Formula = "=M6/30"
$ws.Range("O5:O13").NumberFormat = "0.00%"

# --- "60Sec" table (Q:T): Player 1 win counts ---
$p1Wins60Sec = @(17, 13, 11, 16, 14, 13, 14, 15, 20)
for ($i = 0; $i -lt $p1Wins60Sec.Length; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 18).Value = $p1Wins60Sec[$i]     # column R
}
$ws.Range("S5").Formula = "=30-R5"
$ws.Range("T5").Formula = "=R5/30"
$ws.Range("S6:S13").Formula = "=30-R6"
$ws.Range("T6:T13").Formula = "=R6/30"
$ws.Range("T5:T13").NumberFormat = "0.00%"

# --- Selection / view state to match the saved workbook ---
$null = $ws.Range("U21").Select()
